# Refresh the crypto price/volume snapshot (and restore the correct
# Polkadot / WrappedEther row ordering) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.253.90'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.692.32'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.09'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +13.62%  '
$ws.Range("E9").Value = '  +4.81%  '
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0892'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '1.929.90'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.689.21'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("E15").Value = '  +5.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("D17").Value = '27.244.43'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '238.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  +1.82%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  +3.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.92%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("D32").Value = '1.581.75'
$ws.Range("E32").Value = '  +6.96%  '
$ws.Range("E33").Value = '  +2.87%  '
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("E35").Value = '  +0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.959'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.21%  '
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("E40").Value = '  +4.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.78'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("E44").Value = '  -2.49%  '
$ws.Range("D45").Value = '1.840.79'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '91.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("E48").Value = '  +6.34%  '
$ws.Range("E49").Value = '  +2.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.21%  '
$ws.Range("E51").Value = '  +3.33%  '
